# "Selection of which month to mark attendance of selected"
#
# Sheet1 is an attendance tracker: column A is the person's name and
# columns B..M are Jan..Dec. The edit fills in attendance figures for
# Feb/Mar (C/D) for every person, tweaks the Jan (B) "month selector"
# values, resets the Apr (E) column to 0, and locks the cells that are
# no longer meant to be edited (the now-computed/selected months).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# hansraj (row 2)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("E2").Locked = $true

# riguda (row 3)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("D3").Locked = $true
$ws.Range("E3").Value = 0
$ws.Range("E3").Locked = $true

# chikne (row 4)
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D4").Locked = $true
$ws.Range("E4").Value = 0
$ws.Range("E4").Locked = $true

# ejas (row 5)
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D5").Locked = $true
$ws.Range("E5").Value = 0
$ws.Range("E5").Locked = $true

# rahul (row 6)
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("D6").Locked = $true
$ws.Range("E6").Value = 0
$ws.Range("E6").Locked = $true

# jayesh (row 7)
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E7").Locked = $true
